$d = $word.ActiveDocument

# --- Change 1: drop the trailing period on the "Once you've changed..." paragraph
$p5 = $d.Paragraphs.Item(5)
$null = $p5.Range.Find.Execute("task.", $true, $false, $false, $false, $false, $true, 1, $false, "task", 2)

# --- Changes 2 & 3: turn the "..." paragraph and the empty paragraph that follows it
#     into two fully-written paragraphs about version management. We replace both
#     paragraphs in a single InsertXML call (via a Word "single file XML" package
#     fragment) so the paragraph count stays the same and each paragraph ends up
#     with clean, separate runs (no stray run-formatting markers).
$p6 = $d.Paragraphs.Item(6)
$p7 = $d.Paragraphs.Item(7)
$target = $d.Range($p6.Range.Start, $p7.Range.End)

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">Version management is the process of organizing, </w:t></w:r>
<w:r><w:t>maintaining,</w:t></w:r>
<w:r><w:t xml:space="preserve"> and tracking different versions of software or code. It helps in keeping track of changes made to the software or codebase and helps the team to work collaboratively by allowing them to work on different versions of the same codebase. </w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">Version management systems typically include features such as version control, code branching, merging, </w:t></w:r>
<w:r><w:t>labelling</w:t></w:r>
<w:r><w:t xml:space="preserve">, and release management. These features allow developers to keep track of changes, revert to previous versions if necessary, and work on new features without impacting the existing codebase.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.InsertXML($xml)
